$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "IsAvailable" header with a new "Availability" column header in D1,
# and give it a wrap-text style.
$ws.Range("D1").Value = "Availability"
$ws.Range("D1").WrapText = $true

# Widen column D to fit the new header.
$ws.Columns.Item(4).ColumnWidth = 10.6

# Populate the new Availability column (D2:D13) with 0 for every product row.
$ws.Range("D2:D13").Value = 0

# Move the sheet selection onto the new column.
$ws.Range("D5:D7").Select()
